$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.155.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.774.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.33%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "624.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.61%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.771.22"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("E9").Value = "  +1.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("E11").Value = "  +3.34%  "
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.411.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.774.43"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.119.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.68%  "
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "467.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.57"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("E23").Value = "  +2.23%  "
$ws.Range("E24").Value = "  +3.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("E27").Value = "  +4.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.922.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.91%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.42%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.80"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.168"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +16.12%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.725.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("E38").Value = "  +0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.101"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.60%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.966"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("B45").Value = "Arweave"
$ws.Range("C45").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.70%  "
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.297"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "152.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.70"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("E49").Value = "  +4.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.76%  "
$ws.Range("E51").Value = "  +0.80%  "
